$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 2
$ws.Range("G3").Value2 = 1
$ws.Range("G4").Value2 = 2
$ws.Range("G5").Value2 = 1
$ws.Range("G6").Value2 = 1
$ws.Range("G7").Value2 = 1
$ws.Range("G8").Value2 = 0
$ws.Range("G9").Value2 = 2
$ws.Range("G10").Value2 = 1
$ws.Range("G11").Value2 = 1
$ws.Range("G12").Value2 = 2
$ws.Range("G13").Value2 = 1
$ws.Range("G14").Value2 = 2
$ws.Range("G15").Value2 = 1
$ws.Range("G16").Value2 = 1
$ws.Range("G17").Value2 = 1
$ws.Range("G18").Value2 = 0
$ws.Range("G19").Value2 = 2
$ws.Range("G20").Value2 = 0
$ws.Range("G21").Value2 = 1
$ws.Range("G22").Value2 = 1
$ws.Range("G23").Value2 = 2
$ws.Range("G24").Value2 = 0
$ws.Range("G25").Value2 = 1
$ws.Range("G26").Value2 = 1
$ws.Range("G27").Value2 = 0
$ws.Range("G29").Value2 = 1
$ws.Range("G30").Value2 = 3
$ws.Range("G31").Value2 = 2
$ws.Range("G32").Value2 = 0
$ws.Range("G33").Value2 = 2
$ws.Range("G34").Value2 = 0
$ws.Range("G35").Value2 = 1
$ws.Range("G36").Value2 = 1
$ws.Range("G37").Value2 = 0
$ws.Range("G38").Value2 = 3
$ws.Range("G39").Value2 = 2
$ws.Range("G40").Value2 = 2
$ws.Range("G41").Value2 = 0
$ws.Range("G42").Value2 = 1
$ws.Range("G43").Value2 = 1
$ws.Range("G44").Value2 = 2
$ws.Range("G45").Value2 = 0
$ws.Range("G46").Value2 = 2
$ws.Range("G47").Value2 = 1
$ws.Range("G48").Value2 = 1
$ws.Range("G49").Value2 = 3
$ws.Range("G50").Value2 = 0
$ws.Range("G51").Value2 = 0
$ws.Range("G52").Value2 = 1
$ws.Range("G53").Value2 = 0
$ws.Range("G54").Value2 = 1
$ws.Range("G55").Value2 = 1
$ws.Range("G56").Value2 = 2
$ws.Range("G57").Value2 = 2
$ws.Range("G58").Value2 = 2
$ws.Range("G59").Value2 = 0
$ws.Range("G60").Value2 = 0
$ws.Range("G61").Value2 = 1
$ws.Range("G62").Value2 = 2
$ws.Range("G63").Value2 = 1
$ws.Range("G64").Value2 = 1
$ws.Range("G65").Value2 = 2
$ws.Range("G66").Value2 = 1
$ws.Range("G67").Value2 = 1
$ws.Range("G68").Value2 = 1
$ws.Range("G69").Value2 = 0
$ws.Range("G70").Value2 = 1
$ws.Range("G71").Value2 = 0
$ws.Range("G72").Value2 = 0
$ws.Range("G73").Value2 = 1
$ws.Range("G74").Value2 = 2
$ws.Range("G76").Value2 = 1
$ws.Range("G79").Value2 = 1
$ws.Range("G80").Value2 = 2
